$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 235; this shifts existing rows 235-352 down to 236-353
$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with the new weekly record
$ws.Range("A235").Value = 10
$ws.Range("B235").Value = "Vega Modelo de Temuco"
$ws.Range("C235").Value = "La Araucanía"
$ws.Range("D235").Value = 45001
$ws.Range("E235").Value = 9
$ws.Range("F235").Value = 100112039
$ws.Range("G235").Value = "Ciboulette"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 100
$ws.Range("K235").Value = 5000
$ws.Range("L235").Value = 5000
$ws.Range("M235").Value = 5000
$ws.Range("N235").Value = '$/docena de atados'
$ws.Range("O235").Value = "Provincia de Cautín"
$ws.Range("P235").Value = 1667
$ws.Range("Q235").Value = 3
$ws.Range("R235").Value = "Hortaliza"
